$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Enter new data rows (pieces #5 / row 6, and piece #6 / row 7 partially) ---
$ws.Range("C6").Value = 7.6500000953674316
$ws.Range("D6").Value = 3
$ws.Range("D7").Value = 14

# --- Add the running-total formulas in row 2 ---
$ws.Range("F2").Formula = "=SUM(C2:C300)"
$ws.Range("G2").Formula = "=SUM(D2:D300)"

# Freeze the current calculation results (matches the as-authored workbook,
# where F2/G2 were computed before the C7 time entry was typed in) by
# switching to manual calculation before saving.
$excel.Calculation = -4135
$wb.Save()

# --- Finish entering the last data row ---
$ws.Range("C7").Value = 2.5550000667572021

# --- Move the active selection to match the author's final cursor position ---
$null = $ws.Range("F8").Select()
